# [IMP] z0bug_odoo invoice data
#
# The chart-of-accounts codes referenced by the invoice-line "account_id"
# column (G) used the internal "z0bug.coa_*" ids. Re-point them at the
# external "external.*" ids instead. Also refresh the sheet's remembered
# view state (frozen-pane scroll position + current selection) to match
# where the user left off reviewing column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.Application.ActiveWindow

# --- 1. Update the account_id values (column G) ---------------------------
$replacements = @{
    "z0bug.coa_512000" = "external.512000"
    "z0bug.coa_510200" = "external.510200"
    "z0bug.coa_510000" = "external.510000"
    "z0bug.coa_510100" = "external.510100"
    "z0bug.coa_623460" = "external.623460"
    "z0bug.coa_123380" = "external.123380"
    "z0bug.coa_610100" = "external.610100"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}

# --- 2. Update the frozen-pane scroll position -----------------------------
# Keep the header row frozen but bring row 50 to the top of the
# scrollable area (was row 32).
$win.FreezePanes = $false
$null = $ws.Range("A50").Select()
$win.FreezePanes = $true

# --- 3. Update the remembered selection ------------------------------------
# The reviewer was looking at the account_id column values (G2:G57).
$null = $ws.Range("G2:G57").Select()
